# Auto update Excel log
# Appends new sensor-log rows to the PIR, Humidity, Temperature, Proximity
# and Camera sheets (all for 2026-02-01, ~14:11-14:12 timeframe).
#
# Column A holds date-looking text ("2026-02-01") and (on the Humidity
# sheet) column E holds percentage-looking text ("78.3%"). Excel's normal
# Range.Value assignment auto-coerces those literals into date/percentage
# numbers, but the source log stores them as plain text. Pre-formatting the
# destination cells as Text ("@") before assigning the value keeps them as
# the literal strings the log expects.

$wb = $excel.ActiveWorkbook

function Add-LogRows {
    param(
        [string]$SheetName,
        [object[]]$Rows,
        [bool]$ProtectValueColumn = $false
    )
    # NOTE: callers must pass arguments positionally (not as -SheetName /
    # -Rows / ...). This runtime's param-binder splats an array argument
    # across positional slots when it is bound via a named parameter,
    # which silently corrupts $Rows. Positional calls bind correctly.

    $ws = $wb.Worksheets.Item($SheetName)

    $firstRow = [int]$Rows[0][0]
    $lastRow = [int]$Rows[$Rows.Count - 1][0]

    # Force the Date column to be read back as literal text instead of
    # being auto-parsed into a date serial number.
    $ws.Range("A$firstRow`:A$lastRow").NumberFormat = "@"
    if ($ProtectValueColumn) {
        # Humidity's Value column contains literal "NN.N%" strings that
        # Excel would otherwise coerce into a percentage number.
        $ws.Range("E$firstRow`:E$lastRow").NumberFormat = "@"
    }

    foreach ($r in $Rows) {
        $rowNum = [int]$r[0]
        $ws.Cells.Item($rowNum, 1).Value = $r[1]
        $ws.Cells.Item($rowNum, 2).Value = $r[2]
        $ws.Cells.Item($rowNum, 3).Value = $r[3]
        $ws.Cells.Item($rowNum, 4).Value = $r[4]
        $ws.Cells.Item($rowNum, 5).Value = $r[5]
        $ws.Cells.Item($rowNum, 6).Value = $r[6]
    }
}

# --- PIR sheet: rows 214-226 --------------------------------------------
$pirRows = @(
    @(214, "2026-02-01", "14:11:32", "14:00", "Bathroom", "No Motion", "Inactive"),
    @(215, "2026-02-01", "14:11:34", "14:00", "Bathroom", "No Motion", "Inactive"),
    @(216, "2026-02-01", "14:11:49", "14:00", "Bathroom", "Motion Detected", "Active"),
    @(217, "2026-02-01", "14:11:50", "14:00", "Bathroom", "No Motion", "Inactive"),
    @(218, "2026-02-01", "14:11:52", "14:00", "Bathroom", "Motion Detected", "Active"),
    @(219, "2026-02-01", "14:11:53", "14:00", "Bathroom", "No Motion", "Inactive"),
    @(220, "2026-02-01", "14:11:58", "14:00", "Bathroom", "No Motion", "Inactive"),
    @(221, "2026-02-01", "14:12:03", "14:00", "Bathroom", "No Motion", "Inactive"),
    @(222, "2026-02-01", "14:12:08", "14:00", "Bathroom", "No Motion", "Inactive"),
    @(223, "2026-02-01", "14:12:13", "14:00", "Bathroom", "No Motion", "Inactive"),
    @(224, "2026-02-01", "14:12:19", "14:00", "Bathroom", "No Motion", "Inactive"),
    @(225, "2026-02-01", "14:12:24", "14:00", "Bathroom", "No Motion", "Inactive"),
    @(226, "2026-02-01", "14:12:29", "14:00", "Bathroom", "No Motion", "Inactive")
)
Add-LogRows "PIR" $pirRows

# --- Humidity sheet: rows 134-145 ---------------------------------------
$humidityRows = @(
    @(134, "2026-02-01", "14:11:33", "14:00", "Bathroom", "78.3%", "Active"),
    @(135, "2026-02-01", "14:11:35", "14:00", "Bathroom", "77.3%", "Active"),
    @(136, "2026-02-01", "14:11:49", "14:00", "Bathroom", "78.3%", "Active"),
    @(137, "2026-02-01", "14:11:51", "14:00", "Bathroom", "77.3%", "Active"),
    @(138, "2026-02-01", "14:11:52", "14:00", "Bathroom", "78.3%", "Active"),
    @(139, "2026-02-01", "14:12:01", "14:00", "Bathroom", "78.2%", "Active"),
    @(140, "2026-02-01", "14:12:06", "14:00", "Bathroom", "77.3%", "Active"),
    @(141, "2026-02-01", "14:12:11", "14:00", "Bathroom", "78.3%", "Active"),
    @(142, "2026-02-01", "14:12:16", "14:00", "Bathroom", "77.3%", "Active"),
    @(143, "2026-02-01", "14:12:21", "14:00", "Bathroom", "78.2%", "Active"),
    @(144, "2026-02-01", "14:12:26", "14:00", "Bathroom", "77.2%", "Active"),
    @(145, "2026-02-01", "14:12:31", "14:00", "Bathroom", "78.2%", "Active")
)
Add-LogRows "Humidity" $humidityRows $true

# --- Temperature sheet: rows 55-66 --------------------------------------
$temperatureRows = @(
    @(55, "2026-02-01", "14:11:33", "14:00", "Bathroom", "29.4C", "Active"),
    @(56, "2026-02-01", "14:11:36", "14:00", "Bathroom", "29.3C", "Active"),
    @(57, "2026-02-01", "14:11:50", "14:00", "Bathroom", "29.4C", "Active"),
    @(58, "2026-02-01", "14:11:51", "14:00", "Bathroom", "29.3C", "Active"),
    @(59, "2026-02-01", "14:11:53", "14:00", "Bathroom", "29.4C", "Active"),
    @(60, "2026-02-01", "14:12:01", "14:00", "Bathroom", "29.3C", "Active"),
    @(61, "2026-02-01", "14:12:06", "14:00", "Bathroom", "29.3C", "Active"),
    @(62, "2026-02-01", "14:12:11", "14:00", "Bathroom", "29.4C", "Active"),
    @(63, "2026-02-01", "14:12:16", "14:00", "Bathroom", "29.4C", "Active"),
    @(64, "2026-02-01", "14:12:22", "14:00", "Bathroom", "29.4C", "Active"),
    @(65, "2026-02-01", "14:12:26", "14:00", "Bathroom", "29.4C", "Active"),
    @(66, "2026-02-01", "14:12:31", "14:00", "Bathroom", "29.4C", "Active")
)
Add-LogRows "Temperature" $temperatureRows

# --- Proximity sheet: row 27 --------------------------------------------
# NOTE: the leading "," forces this to stay an array-of-one-array. Without
# it, a single-element @() collapses and $proximityRows would become the
# *inner* 7-item row array instead of a 1-item array of rows.
$proximityRows = @(
    ,@(27, "2026-02-01", "14:11:36", "14:00", "Living Room Main Door", "ENTER", "User ENTERED Living Room Main Door")
)
Add-LogRows "Proximity" $proximityRows

# --- Camera sheet: row 16 ------------------------------------------------
$cameraRows = @(
    ,@(16, "2026-02-01", "14:11:48", "14:00", "Living Room Main Door", "Image Captured", "Active")
)
Add-LogRows "Camera" $cameraRows
